# Apache POI wrote new row data into Sheet2 (row 2), replacing the previous
# record (Gary / Torphy / 66S6O@mailinator.com / 9826031438) with a freshly
# "formulated" one read from the workbook model.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A2").Value = "Clemente"
$ws.Range("B2").Value = "Pacocha"
$ws.Range("C2").Value = "Xe4Eu@mailinator.com"
# Mobile numbers are stored as text (leading apostrophe keeps Excel from
# reinterpreting the digit string as a numeric value).
$ws.Range("D2").Value = "'9826098823"
